# Auto-generated script applying scheduled market-data refresh to Jenova_Profits sheets
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("N3").ClearContents()
$ws.Range("H5").Value = 63.5
$ws.Range("I5").Value = 66.583336
$ws.Range("K5").Value = 66.583336
$ws.Range("M5").Value = 48.416664
$ws.Range("H100").Value = 8917.308000000001
$ws.Range("I100").Value = 3865.75
$ws.Range("J100").Value = 16999.8
$ws.Range("K100").Value = 3865.75
$ws.Range("L100").Value = 16999.8
$ws.Range("M100").Value = -3324.75
$ws.Range("N100").Value = -18081.8
$ws.Range("H102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("N102").ClearContents()
$ws.Range("H138").Value = 10676.192
$ws.Range("J138").Value = 10777.648
$ws.Range("L138").Value = 32332.944
$ws.Range("N138").Value = -42612.944

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").ClearContents()
$ws.Range("H32").Value = 20957.195
$ws.Range("I32").Value = 14625.521
$ws.Range("J32").Value = 27288.87
$ws.Range("K32").Value = 14625.521
$ws.Range("L32").Value = 27288.87
$ws.Range("M32").Value = -14338.521
$ws.Range("N32").Value = -27862.87
$ws.Range("H61").Value = 7971.357
$ws.Range("I61").Value = 5859.95
$ws.Range("J61").Value = 13249.875
$ws.Range("K61").Value = 5859.95
$ws.Range("L61").Value = 13249.875
$ws.Range("M61").Value = -5647.95
$ws.Range("N61").Value = -13673.875
$ws.Range("H74").Value = 4552.5713
$ws.Range("I74").Value = 3644.6667
$ws.Range("J74").Value = 10000
$ws.Range("K74").Value = 3644.6667
$ws.Range("L74").Value = 10000
$ws.Range("M74").Value = -2770.6667
$ws.Range("N74").Value = -11748
$ws.Range("H77").Value = 4552.5713
$ws.Range("I77").Value = 3644.6667
$ws.Range("J77").Value = 10000
$ws.Range("K77").Value = 18223.3335
$ws.Range("L77").Value = 50000
$ws.Range("M77").Value = -13855.3335
$ws.Range("N77").Value = -58736
$ws.Range("H80").Value = 86000
$ws.Range("J80").Value = 86000
$ws.Range("L80").Value = 86000
$ws.Range("N80").Value = -87996
$ws.Range("H83").Value = 86000
$ws.Range("J83").Value = 86000
$ws.Range("L83").Value = 258000
$ws.Range("N83").Value = -267984
$ws.Range("H97").Value = 2520
$ws.Range("I97").Value = 2520
$ws.Range("K97").Value = 2520
$ws.Range("M97").Value = -2024
$ws.Range("H102").Value = 3118.0588
$ws.Range("I102").Value = 3013.7334
$ws.Range("J102").Value = 3900.5
$ws.Range("K102").Value = 3013.7334
$ws.Range("L102").Value = 3900.5
$ws.Range("M102").Value = -1391.7334
$ws.Range("N102").Value = -7144.5
$ws.Range("H110").Value = 180706.42
$ws.Range("I110").Value = 201807.2
$ws.Range("K110").Value = 201807.2
$ws.Range("M110").Value = -199762.2
$ws.Range("H122").Value = 3770.077
$ws.Range("I122").Value = 2302.75
$ws.Range("K122").Value = 6908.25
$ws.Range("M122").Value = -4458.25
$ws.Range("H132").Value = 4359.3877
$ws.Range("I132").Value = 3545.3333
$ws.Range("J132").Value = 8197.071
$ws.Range("K132").Value = 10635.9999
$ws.Range("L132").Value = 24591.213
$ws.Range("M132").Value = -8105.999899999999
$ws.Range("N132").Value = -29651.213
$ws.Range("H136").Value = 7971.357
$ws.Range("I136").Value = 5859.95
$ws.Range("J136").Value = 13249.875
$ws.Range("K136").Value = 17579.85
$ws.Range("L136").Value = 39749.625
$ws.Range("M136").Value = -15029.85
$ws.Range("N136").Value = -44849.625

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1548182.2
$ws.Range("J86").Value = 3498.5
$ws.Range("L86").Value = 3498.5
$ws.Range("N86").Value = -5744.5
$ws.Range("H89").Value = 1548182.2
$ws.Range("J89").Value = 3498.5
$ws.Range("L89").Value = 17492.5
$ws.Range("N89").Value = -28724.5
$ws.Range("H99").Value = 1236.2667
$ws.Range("I99").Value = 1195.7693
$ws.Range("K99").Value = 1195.7693
$ws.Range("M99").Value = 302.2307000000001
$ws.Range("H134").Value = 20603.46
$ws.Range("I134").Value = 4833.852
$ws.Range("K134").Value = 14501.556
$ws.Range("M134").Value = -11966.556
$ws.Range("H139").Value = 120000
$ws.Range("J139").Value = 120000
$ws.Range("L139").Value = 120000
$ws.Range("N139").Value = -130280

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 26815.79
$ws.Range("I6").Value = 29071.428
$ws.Range("K6").Value = 29071.428
$ws.Range("M6").Value = -28958.428
$ws.Range("H7").Value = 563
$ws.Range("I7").Value = 545.6667
$ws.Range("J7").Value = 615
$ws.Range("K7").Value = 545.6667
$ws.Range("L7").Value = 615
$ws.Range("M7").Value = -432.6667
$ws.Range("N7").Value = -841
$ws.Range("H31").Value = 51581.816
$ws.Range("I31").Value = 4249.6875
$ws.Range("J31").Value = 177800.83
$ws.Range("K31").Value = 4249.6875
$ws.Range("L31").Value = 177800.83
$ws.Range("M31").Value = -3954.6875
$ws.Range("N31").Value = -178390.83
$ws.Range("H34").Value = 51581.816
$ws.Range("I34").Value = 4249.6875
$ws.Range("J34").Value = 177800.83
$ws.Range("K34").Value = 4249.6875
$ws.Range("L34").Value = 177800.83
$ws.Range("M34").Value = -4047.6875
$ws.Range("N34").Value = -178204.83
$ws.Range("H50").Value = 22939.5
$ws.Range("J50").Value = 56000
$ws.Range("L50").Value = 56000
$ws.Range("N50").Value = -57250

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 1252.75
$ws.Range("I122").Value = 1003.6667
$ws.Range("K122").Value = 9033.0003
$ws.Range("M122").Value = -6583.0003

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 597697.5600000001
$ws.Range("I113").Value = 1253238.6
$ws.Range("J113").Value = 14994.444
$ws.Range("K113").Value = 1253238.6
$ws.Range("L113").Value = 14994.444
$ws.Range("M113").Value = -1251068.6
$ws.Range("N113").Value = -19334.444
$ws.Range("H122").Value = 4299.2856
$ws.Range("I122").Value = 3599
$ws.Range("J122").Value = 6050
$ws.Range("K122").Value = 10797
$ws.Range("L122").Value = 18150
$ws.Range("M122").Value = -8347
$ws.Range("N122").Value = -23050
$ws.Range("H132").Value = 50125.156
$ws.Range("I132").Value = 5586.225
$ws.Range("J132").Value = 406436.6
$ws.Range("K132").Value = 16758.675
$ws.Range("L132").Value = 1219309.8
$ws.Range("M132").Value = -14228.675
$ws.Range("N132").Value = -1224369.8
$ws.Range("H137").Value = 51334.625
$ws.Range("J137").Value = 49996.5
$ws.Range("L137").Value = 49996.5
$ws.Range("N137").Value = -60196.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("I82").Value = 1050
$ws.Range("K82").Value = 1050
$ws.Range("M82").Value = -689
$ws.Range("I85").Value = 1050
$ws.Range("K85").Value = 1050
$ws.Range("M85").Value = 198

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 71439560
$ws.Range("I122").Value = 166680820
$ws.Range("K122").Value = 500042460
$ws.Range("M122").Value = -500040010
$ws.Range("H132").Value = 27917.977
$ws.Range("I132").Value = 4408.9355
$ws.Range("J132").Value = 88649.664
$ws.Range("K132").Value = 13226.8065
$ws.Range("L132").Value = 265948.992
$ws.Range("M132").Value = -10696.8065
$ws.Range("N132").Value = -271008.992
$ws.Range("H136").Value = 268142.38
$ws.Range("I136").Value = 351171.44
$ws.Range("J136").Value = 126504.53
$ws.Range("K136").Value = 1053514.32
$ws.Range("L136").Value = 379513.59
$ws.Range("M136").Value = -1050964.32
$ws.Range("N136").Value = -384613.59
